$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5491146666666666
$ws.Range("H2").Value = 1.647344
$ws.Range("I2").Value = 0.006868658684314029
$ws.Range("J2").Value = 0.006868658684314029
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09834766666666667
$ws.Range("N2").Value = 0.295043
$ws.Range("O2").Value = 0.2818566198948398
$ws.Range("P2").Value = 0.2818566198948398
$ws.Range("Q2").Value = 0.05400414619911111
$ws.Range("R2").Value = 0.486037315792
$ws.Range("S2").Value = 0.00193597691997209
$ws.Range("T2").Value = 0.00193597691997209

$ws.Range("G3").Value = 0.5491146666666666
$ws.Range("H3").Value = 1.647344
$ws.Range("I3").Value = 0.006868658684314029
$ws.Range("J3").Value = 0.006868658684314029
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2505803333333333
$ws.Range("N3").Value = 0.751741
$ws.Range("O3").Value = 0.7181433801051602
$ws.Range("P3").Value = 0.7181433801051602
$ws.Range("Q3").Value = 0.1375973362115556
$ws.Range("R3").Value = 1.238376025904
$ws.Range("S3").Value = 0.004932681764341939
$ws.Range("T3").Value = 0.004932681764341939

$ws.Range("I4").Value = 0.4602735740408745
$ws.Range("J4").Value = 0.4602735740408745
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09834766666666667
$ws.Range("N4").Value = 0.295043
$ws.Range("O4").Value = 0.2818566198948398
$ws.Range("P4").Value = 0.2818566198948398
$ws.Range("Q4").Value = 3.618855227274
$ws.Range("R4").Value = 32.569697045466
$ws.Range("S4").Value = 0.1297311538060782
$ws.Range("T4").Value = 0.1297311538060782

$ws.Range("I5").Value = 0.4602735740408745
$ws.Range("J5").Value = 0.4602735740408745
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2505803333333333
$ws.Range("N5").Value = 0.751741
$ws.Range("O5").Value = 0.7181433801051602
$ws.Range("P5").Value = 0.7181433801051602
$ws.Range("Q5").Value = 9.220492766838001
$ws.Range("R5").Value = 82.984434901542
$ws.Range("S5").Value = 0.3305424202347964
$ws.Range("T5").Value = 0.3305424202347964

$ws.Range("G6").Value = 19.686315
$ws.Range("H6").Value = 59.058945
$ws.Range("I6").Value = 0.2462483461017703
$ws.Range("J6").Value = 0.2462483461017702
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.09834766666666667
$ws.Range("N6").Value = 0.295043
$ws.Range("O6").Value = 0.2818566198948398
$ws.Range("P6").Value = 0.2818566198948398
$ws.Range("Q6").Value = 1.936103145515
$ws.Range("R6").Value = 17.424928309635
$ws.Range("S6").Value = 0.06940672648693962
$ws.Range("T6").Value = 0.06940672648693962

$ws.Range("G7").Value = 19.686315
$ws.Range("H7").Value = 59.058945
$ws.Range("I7").Value = 0.2462483461017703
$ws.Range("J7").Value = 0.2462483461017702
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2505803333333333
$ws.Range("N7").Value = 0.751741
$ws.Range("O7").Value = 0.7181433801051602
$ws.Range("P7").Value = 0.7181433801051602
$ws.Range("Q7").Value = 4.933003374805001
$ws.Range("R7").Value = 44.397030373245
$ws.Range("S7").Value = 0.1768416196148307
$ws.Range("T7").Value = 0.1768416196148306

$ws.Range("G8").Value = 22.91297966666667
$ws.Range("H8").Value = 68.738939
$ws.Range("I8").Value = 0.2866094211730412
$ws.Range("J8").Value = 0.2866094211730412
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.09834766666666667
$ws.Range("N8").Value = 0.295043
$ws.Range("O8").Value = 0.2818566198948398
$ws.Range("P8").Value = 0.2818566198948398
$ws.Range("Q8").Value = 2.253438086597445
$ws.Range("R8").Value = 20.280942779377
$ws.Range("S8").Value = 0.08078276268184992
$ws.Range("T8").Value = 0.08078276268184993

$ws.Range("G9").Value = 22.91297966666667
$ws.Range("H9").Value = 68.738939
$ws.Range("I9").Value = 0.2866094211730412
$ws.Range("J9").Value = 0.2866094211730412
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2505803333333333
$ws.Range("N9").Value = 0.751741
$ws.Range("O9").Value = 0.7181433801051602
$ws.Range("P9").Value = 0.7181433801051602
$ws.Range("Q9").Value = 5.741542082533223
$ws.Range("R9").Value = 51.673878742799
$ws.Range("S9").Value = 0.2058266584911913
$ws.Range("T9").Value = 0.2058266584911913

